$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string that begins with a leading apostrophe.
# Excel's normal cell-input parsing treats a leading "'" as the quote-prefix
# ("format as text") marker and strips it from the stored value, which also
# silently stamps the cell with a quotePrefix style. To store the apostrophe
# as literal text (matching the target data) we instead build the text with
# a formula (CHAR(39) = an apostrophe character) and then convert that
# formula to a plain value in place via copy / paste-special-values - this
# keeps the literal leading apostrophe without leaving a formula behind or
# picking up the quote-prefix cell style.
function Set-LiteralText {
    param([string]$addr, $value)
    $ws.Range($addr).Value = $value
}

function Set-LeadingApostropheText {
    param([string]$addr, [string]$afterFirstQuote)
    $ws.Range($addr).Formula = '=CHAR(39)&"' + $afterFirstQuote + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# Header row
Set-LiteralText "A1" "Review"
Set-LiteralText "B1" "name"
Set-LiteralText "C1" "date"
Set-LiteralText "D1" "length"
Set-LiteralText "E1" "policy"
Set-LiteralText "F1" "policyType"

# Row 2
Set-LeadingApostropheText "A2" "it just works'"
Set-LeadingApostropheText "B2" "steve'"
Set-LeadingApostropheText "C2" "9/10/21'"
Set-LiteralText "D2" 3
Set-LeadingApostropheText "E2" "Vehicle'"
Set-LiteralText "F2" "Car'"

# Row 3
Set-LeadingApostropheText "A3" "fast service after a freak boating accident'"
Set-LiteralText "B3" " 'gonzalez'"
Set-LiteralText "C3" " '3/21/17'"
Set-LiteralText "D3" 12
Set-LiteralText "E3" "Vehicle'"
Set-LiteralText "F3" "Boat'"

$excel.CutCopyMode = $false

# Update selection to match the commit (D6 selected)
$ws.Range("D6").Select() | Out-Null
